$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fields ---
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 08.11.2023"

# --- Row 6 ---
$ws.Range("B6").Value = "10.11."
$ws.Range("C6").Value = "11.11."
$ws.Range("D6").Value = "KARTENZ./10.11 ALDI SUED RO"
$ws.Range("E6").Value = "92,84-"

# --- Row 7 ---
$ws.Range("B7").Value = "13.11."
$ws.Range("C7").Value = "14.11."
$ws.Range("D7").Value = "AMAZON.DE MKTPLC EU BNCTQH"
$ws.Range("E7").Value = "225,20-"

# --- Row 8 ---
$ws.Range("B8").Value = "14.11."
$ws.Range("C8").Value = "15.11."
$ws.Range("D8").Value = "MCDONALDS Lübben"
$ws.Range("E8").Value = "41,08-"

# --- Row 9 ---
$ws.Range("B9").Value = "15.11."
$ws.Range("C9").Value = "16.11."
$ws.Range("D9").Value = "EBAY MKTPLC EU CBPBEY"
$ws.Range("E9").Value = "130,32-"

# --- Row 10: was blank, now becomes a new transaction row ---
# Copy E9's formatting (right-aligned, non-wrapping) onto E10 first,
# matching how the other populated rows are styled.
$ws.Range("E9").Copy()
$ws.Range("E10").PasteSpecial(-4122)

$ws.Range("B10").Value = "18.11."
$ws.Range("C10").Value = "19.11."
$ws.Range("D10").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E10").Value = "25,08-"

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 21.11.2023"
$ws.Range("E12").Value = "514,52-"

# --- Next billing date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 01.12.2023"
